$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169, shifting existing rows 169:269 down to 170:270
$ws.Rows(169).Insert()

# Fill the new row 169 with the new data record
$ws.Range("A169").Value = 11
$ws.Range("B169").Value = "Vega Monumental Concepción"
$ws.Range("C169").Value = "Bíobío"
$ws.Range("D169").Value = 45176
$ws.Range("E169").Value = 8
$ws.Range("F169").Value = 100112032
$ws.Range("G169").Value = "Zapallo italiano"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 100
$ws.Range("K169").Value = 16000
$ws.Range("L169").Value = 17000
$ws.Range("M169").Value = 16500
$ws.Range("N169").Value = "$/caja 50 unidades"
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 330
$ws.Range("Q169").Value = 50
$ws.Range("R169").Value = "Hortaliza"
